# Auto-generated Excel COM-interop script applying the Aegis_Profits diff
$wb = $excel.ActiveWorkbook

$cellUpdates = @(
    @{Sheet="ALC"; Cell="H100"; Value=1108.0476},
    @{Sheet="ALC"; Cell="I100"; Value=655.9286},
    @{Sheet="ALC"; Cell="K100"; Value=655.9286},
    @{Sheet="ALC"; Cell="M100"; Value=-114.9286},
    @{Sheet="ALC"; Cell="H112"; Value=1316.8235},
    @{Sheet="ALC"; Cell="J112"; Value=1341.5758},
    @{Sheet="ALC"; Cell="L112"; Value=4024.7274},
    @{Sheet="ALC"; Cell="N112"; Value=-6240.7274},
    @{Sheet="ALC"; Cell="H116"; Value=12444.333},
    @{Sheet="ALC"; Cell="I116"; Value=12444.333},
    @{Sheet="ALC"; Cell="K116"; Value=12444.333},
    @{Sheet="ALC"; Cell="M116"; Value=-9002.333000000001},
    @{Sheet="ALC"; Cell="H129"; Value=2362.9849},
    @{Sheet="ALC"; Cell="I129"; Value=7211.6665},
    @{Sheet="ALC"; Cell="J129"; Value=936.902},
    @{Sheet="ALC"; Cell="K129"; Value=21634.9995},
    @{Sheet="ALC"; Cell="L129"; Value=2810.706},
    @{Sheet="ALC"; Cell="M129"; Value=-16634.9995},
    @{Sheet="ALC"; Cell="N129"; Value=-12810.706},
    @{Sheet="ALC"; Cell="H135"; Value=1366.1459},
    @{Sheet="ALC"; Cell="I135"; Value=569.1875},
    @{Sheet="ALC"; Cell="K135"; Value=5122.6875},
    @{Sheet="ALC"; Cell="M135"; Value=-2587.6875},
    @{Sheet="ALC"; Cell="H138"; Value=3356.0942},
    @{Sheet="ALC"; Cell="I138"; Value=2613.1428},
    @{Sheet="ALC"; Cell="J138"; Value=3622.795},
    @{Sheet="ALC"; Cell="K138"; Value=7839.428400000001},
    @{Sheet="ALC"; Cell="L138"; Value=10868.385},
    @{Sheet="ALC"; Cell="M138"; Value=-2699.428400000001},
    @{Sheet="ALC"; Cell="N138"; Value=-21148.385},
    @{Sheet="ARM"; Cell="H32"; Value=41057.1},
    @{Sheet="ARM"; Cell="I32"; Value=11112.488},
    @{Sheet="ARM"; Cell="J32"; Value=225002.58},
    @{Sheet="ARM"; Cell="K32"; Value=11112.488},
    @{Sheet="ARM"; Cell="L32"; Value=225002.58},
    @{Sheet="ARM"; Cell="M32"; Value=-10825.488},
    @{Sheet="ARM"; Cell="N32"; Value=-225576.58},
    @{Sheet="ARM"; Cell="H112"; Value=14633.333},
    @{Sheet="ARM"; Cell="J112"; Value=14633.333},
    @{Sheet="ARM"; Cell="L112"; Value=14633.333},
    @{Sheet="ARM"; Cell="N112"; Value=-17587.333},
    @{Sheet="ARM"; Cell="H114"; Value=27525},
    @{Sheet="ARM"; Cell="J114"; Value=27525},
    @{Sheet="ARM"; Cell="L114"; Value=27525},
    @{Sheet="ARM"; Cell="N114"; Value=-36203},
    @{Sheet="BSM"; Cell="H64"; Value=1866.6666},
    @{Sheet="BSM"; Cell="I64"; Value=700},
    @{Sheet="BSM"; Cell="J64"; Value=2046.1538},
    @{Sheet="BSM"; Cell="K64"; Value=700},
    @{Sheet="BSM"; Cell="L64"; Value=2046.1538},
    @{Sheet="BSM"; Cell="M64"; Value=-475},
    @{Sheet="BSM"; Cell="N64"; Value=-2496.1538},
    @{Sheet="BSM"; Cell="H67"; Value=1866.6666},
    @{Sheet="BSM"; Cell="I67"; Value=700},
    @{Sheet="BSM"; Cell="J67"; Value=2046.1538},
    @{Sheet="BSM"; Cell="K67"; Value=700},
    @{Sheet="BSM"; Cell="L67"; Value=2046.1538},
    @{Sheet="BSM"; Cell="M67"; Value=80},
    @{Sheet="BSM"; Cell="N67"; Value=-3606.1538},
    @{Sheet="BSM"; Cell="H68"; Value=40000},
    @{Sheet="BSM"; Cell="J68"; Value=40000},
    @{Sheet="BSM"; Cell="L68"; Value=40000},
    @{Sheet="BSM"; Cell="N68"; Value=-41622},
    @{Sheet="BSM"; Cell="H71"; Value=40000},
    @{Sheet="BSM"; Cell="J71"; Value=40000},
    @{Sheet="BSM"; Cell="L71"; Value=120000},
    @{Sheet="BSM"; Cell="N71"; Value=-128112},
    @{Sheet="BSM"; Cell="H134"; Value=17932.354},
    @{Sheet="BSM"; Cell="I134"; Value=20127.828},
    @{Sheet="BSM"; Cell="K134"; Value=60383.484},
    @{Sheet="BSM"; Cell="M134"; Value=-57848.484},
    @{Sheet="CRP"; Cell="H22"; Value=324.35294},
    @{Sheet="CRP"; Cell="I22"; Value=303.25},
    @{Sheet="CRP"; Cell="J22"; Value=375},
    @{Sheet="CRP"; Cell="K22"; Value=303.25},
    @{Sheet="CRP"; Cell="L22"; Value=375},
    @{Sheet="CRP"; Cell="M22"; Value=46.75},
    @{Sheet="CRP"; Cell="N22"; Value=-1075},
    @{Sheet="CRP"; Cell="H31"; Value=39359.344},
    @{Sheet="CRP"; Cell="I31"; Value=640.65},
    @{Sheet="CRP"; Cell="J31"; Value=82380.11},
    @{Sheet="CRP"; Cell="K31"; Value=640.65},
    @{Sheet="CRP"; Cell="L31"; Value=82380.11},
    @{Sheet="CRP"; Cell="M31"; Value=-345.65},
    @{Sheet="CRP"; Cell="N31"; Value=-82970.11},
    @{Sheet="CRP"; Cell="H34"; Value=39359.344},
    @{Sheet="CRP"; Cell="I34"; Value=640.65},
    @{Sheet="CRP"; Cell="J34"; Value=82380.11},
    @{Sheet="CRP"; Cell="K34"; Value=640.65},
    @{Sheet="CRP"; Cell="L34"; Value=82380.11},
    @{Sheet="CRP"; Cell="M34"; Value=-438.65},
    @{Sheet="CRP"; Cell="N34"; Value=-82784.11},
    @{Sheet="CRP"; Cell="H62"; Value=2516.6667},
    @{Sheet="CRP"; Cell="I62"; Value=2200},
    @{Sheet="CRP"; Cell="J62"; Value=2545.4546},
    @{Sheet="CRP"; Cell="K62"; Value=2200},
    @{Sheet="CRP"; Cell="L62"; Value=2545.4546},
    @{Sheet="CRP"; Cell="M62"; Value=-1576},
    @{Sheet="CRP"; Cell="N62"; Value=-3793.4546},
    @{Sheet="CRP"; Cell="H65"; Value=2516.6667},
    @{Sheet="CRP"; Cell="I65"; Value=2200},
    @{Sheet="CRP"; Cell="J65"; Value=2545.4546},
    @{Sheet="CRP"; Cell="K65"; Value=11000},
    @{Sheet="CRP"; Cell="L65"; Value=12727.273},
    @{Sheet="CRP"; Cell="M65"; Value=-7880},
    @{Sheet="CRP"; Cell="N65"; Value=-18967.273},
    @{Sheet="CRP"; Cell="H80"; Value=9629.6},
    @{Sheet="CRP"; Cell="J80"; Value=9629.6},
    @{Sheet="CRP"; Cell="L80"; Value=9629.6},
    @{Sheet="CRP"; Cell="N80"; Value=-11875.6},
    @{Sheet="CRP"; Cell="H83"; Value=9629.6},
    @{Sheet="CRP"; Cell="J83"; Value=9629.6},
    @{Sheet="CRP"; Cell="L83"; Value=28888.8},
    @{Sheet="CRP"; Cell="N83"; Value=-40120.8},
    @{Sheet="CRP"; Cell="H122"; Value=1147.2941},
    @{Sheet="CRP"; Cell="I122"; Value=1134.8667},
    @{Sheet="CRP"; Cell="J122"; Value=1240.5},
    @{Sheet="CRP"; Cell="K122"; Value=3404.6001},
    @{Sheet="CRP"; Cell="L122"; Value=3721.5},
    @{Sheet="CRP"; Cell="M122"; Value=-954.6001000000001},
    @{Sheet="CRP"; Cell="N122"; Value=-8621.5},
    @{Sheet="CUL"; Cell="H64"; Value=1878},
    @{Sheet="CUL"; Cell="I64"; Value=756},
    @{Sheet="CUL"; Cell="J64"; Value=3000},
    @{Sheet="CUL"; Cell="K64"; Value=2268},
    @{Sheet="CUL"; Cell="L64"; Value=9000},
    @{Sheet="CUL"; Cell="M64"; Value=-1998},
    @{Sheet="CUL"; Cell="N64"; Value=-9540},
    @{Sheet="CUL"; Cell="H67"; Value=1878},
    @{Sheet="CUL"; Cell="I67"; Value=756},
    @{Sheet="CUL"; Cell="J67"; Value=3000},
    @{Sheet="CUL"; Cell="K67"; Value=2268},
    @{Sheet="CUL"; Cell="L67"; Value=9000},
    @{Sheet="CUL"; Cell="M67"; Value=-1332},
    @{Sheet="CUL"; Cell="N67"; Value=-10872},
    @{Sheet="CUL"; Cell="H132"; Value=346429.53},
    @{Sheet="CUL"; Cell="I132"; Value=515},
    @{Sheet="CUL"; Cell="J132"; Value=627485.0600000001},
    @{Sheet="CUL"; Cell="K132"; Value=4635},
    @{Sheet="CUL"; Cell="L132"; Value=5647365.540000001},
    @{Sheet="CUL"; Cell="M132"; Value=-2105},
    @{Sheet="CUL"; Cell="N132"; Value=-5652425.540000001},
    @{Sheet="GSM"; Cell="H69"; Value=58375},
    @{Sheet="GSM"; Cell="J69"; Value=58375},
    @{Sheet="GSM"; Cell="L69"; Value=58375},
    @{Sheet="GSM"; Cell="N69"; Value=-59873},
    @{Sheet="GSM"; Cell="H72"; Value=58375},
    @{Sheet="GSM"; Cell="J72"; Value=58375},
    @{Sheet="GSM"; Cell="L72"; Value=175125},
    @{Sheet="GSM"; Cell="N72"; Value=-182613},
    @{Sheet="GSM"; Cell="H97"; Value=47621464},
    @{Sheet="GSM"; Cell="I97"; Value=76925870},
    @{Sheet="GSM"; Cell="J97"; Value=1802.75},
    @{Sheet="GSM"; Cell="K97"; Value=76925870},
    @{Sheet="GSM"; Cell="L97"; Value=1802.75},
    @{Sheet="GSM"; Cell="M97"; Value=-76925374},
    @{Sheet="GSM"; Cell="N97"; Value=-2794.75},
    @{Sheet="GSM"; Cell="H103"; Value=54800},
    @{Sheet="GSM"; Cell="J103"; Value=54800},
    @{Sheet="GSM"; Cell="L103"; Value=54800},
    @{Sheet="GSM"; Cell="N103"; Value=-57144},
    @{Sheet="GSM"; Cell="H111"; Value=38902.168},
    @{Sheet="GSM"; Cell="J111"; Value=38902.168},
    @{Sheet="GSM"; Cell="L111"; Value=38902.168},
    @{Sheet="GSM"; Cell="N111"; Value=-45036.168},
    @{Sheet="GSM"; Cell="H112"; Value=40599},
    @{Sheet="GSM"; Cell="J112"; Value=40599},
    @{Sheet="GSM"; Cell="L112"; Value=40599},
    @{Sheet="GSM"; Cell="N112"; Value=-42815},
    @{Sheet="GSM"; Cell="H122"; Value=2404.8823},
    @{Sheet="GSM"; Cell="I122"; Value=1951.3077},
    @{Sheet="GSM"; Cell="J122"; Value=3879},
    @{Sheet="GSM"; Cell="K122"; Value=5853.9231},
    @{Sheet="GSM"; Cell="L122"; Value=11637},
    @{Sheet="GSM"; Cell="M122"; Value=-3403.9231},
    @{Sheet="GSM"; Cell="N122"; Value=-16537},
    @{Sheet="GSM"; Cell="H126"; Value=2602.2693},
    @{Sheet="GSM"; Cell="I126"; Value=2585.25},
    @{Sheet="GSM"; Cell="J126"; Value=2659},
    @{Sheet="GSM"; Cell="K126"; Value=7755.75},
    @{Sheet="GSM"; Cell="L126"; Value=7977},
    @{Sheet="GSM"; Cell="M126"; Value=-5285.75},
    @{Sheet="GSM"; Cell="N126"; Value=-12917},
    @{Sheet="LTW"; Cell="H7"; Value=2585.8125},
    @{Sheet="LTW"; Cell="I7"; Value=1669.8182},
    @{Sheet="LTW"; Cell="J7"; Value=4601},
    @{Sheet="LTW"; Cell="K7"; Value=1669.8182},
    @{Sheet="LTW"; Cell="L7"; Value=4601},
    @{Sheet="LTW"; Cell="M7"; Value=-1557.8182},
    @{Sheet="LTW"; Cell="N7"; Value=-4825},
    @{Sheet="LTW"; Cell="H42"; Value=12559},
    @{Sheet="LTW"; Cell="J42"; Value=12559},
    @{Sheet="LTW"; Cell="L42"; Value=12559},
    @{Sheet="LTW"; Cell="N42"; Value=-13685},
    @{Sheet="LTW"; Cell="H49"; Value=12559},
    @{Sheet="LTW"; Cell="J49"; Value=12559},
    @{Sheet="LTW"; Cell="L49"; Value=12559},
    @{Sheet="LTW"; Cell="N49"; Value=-12853},
    @{Sheet="LTW"; Cell="H59"; Value=0},
    @{Sheet="LTW"; Cell="J59"; Value=0},
    @{Sheet="LTW"; Cell="L59"; Value=0},
    @{Sheet="LTW"; Cell="H69"; Value=36100},
    @{Sheet="LTW"; Cell="J69"; Value=36100},
    @{Sheet="LTW"; Cell="L69"; Value=36100},
    @{Sheet="LTW"; Cell="N69"; Value=-37722},
    @{Sheet="LTW"; Cell="H72"; Value=36100},
    @{Sheet="LTW"; Cell="J72"; Value=36100},
    @{Sheet="LTW"; Cell="L72"; Value=108300},
    @{Sheet="LTW"; Cell="N72"; Value=-116412},
    @{Sheet="LTW"; Cell="H100"; Value=2172.125},
    @{Sheet="LTW"; Cell="I100"; Value=1900},
    @{Sheet="LTW"; Cell="J100"; Value=2335.4},
    @{Sheet="LTW"; Cell="K100"; Value=1900},
    @{Sheet="LTW"; Cell="L100"; Value=2335.4},
    @{Sheet="LTW"; Cell="M100"; Value=-1359},
    @{Sheet="LTW"; Cell="N100"; Value=-3417.4},
    @{Sheet="LTW"; Cell="H110"; Value=28725},
    @{Sheet="LTW"; Cell="J110"; Value=28725},
    @{Sheet="LTW"; Cell="L110"; Value=28725},
    @{Sheet="LTW"; Cell="N110"; Value=-36905},
    @{Sheet="LTW"; Cell="H126"; Value=2585.8125},
    @{Sheet="LTW"; Cell="I126"; Value=1669.8182},
    @{Sheet="LTW"; Cell="J126"; Value=4601},
    @{Sheet="LTW"; Cell="K126"; Value=5009.4546},
    @{Sheet="LTW"; Cell="L126"; Value=13803},
    @{Sheet="LTW"; Cell="M126"; Value=-2539.4546},
    @{Sheet="LTW"; Cell="N126"; Value=-18743},
    @{Sheet="LTW"; Cell="H132"; Value=4059.9},
    @{Sheet="LTW"; Cell="I132"; Value=4383.222},
    @{Sheet="LTW"; Cell="J132"; Value=1150},
    @{Sheet="LTW"; Cell="K132"; Value=13149.666},
    @{Sheet="LTW"; Cell="L132"; Value=3450},
    @{Sheet="LTW"; Cell="M132"; Value=-10619.666},
    @{Sheet="LTW"; Cell="N132"; Value=-8510},
    @{Sheet="WVR"; Cell="H100"; Value=63990.875},
    @{Sheet="WVR"; Cell="I100"; Value=84662.836},
    @{Sheet="WVR"; Cell="J100"; Value=1975},
    @{Sheet="WVR"; Cell="K100"; Value=169325.672},
    @{Sheet="WVR"; Cell="L100"; Value=3950},
    @{Sheet="WVR"; Cell="M100"; Value=-168784.672},
    @{Sheet="WVR"; Cell="N100"; Value=-5032},
    @{Sheet="WVR"; Cell="H107"; Value=77360.766},
    @{Sheet="WVR"; Cell="I107"; Value=453.63635},
    @{Sheet="WVR"; Cell="K107"; Value=1360.90905},
    @{Sheet="WVR"; Cell="M107"; Value=559.09095},
    @{Sheet="WVR"; Cell="H114"; Value=30000},
    @{Sheet="WVR"; Cell="J114"; Value=30000},
    @{Sheet="WVR"; Cell="L114"; Value=30000},
    @{Sheet="WVR"; Cell="N114"; Value=-38678},
    @{Sheet="WVR"; Cell="H132"; Value=4840.442},
    @{Sheet="WVR"; Cell="I132"; Value=3735.35},
    @{Sheet="WVR"; Cell="J132"; Value=5801.391},
    @{Sheet="WVR"; Cell="K132"; Value=11206.05},
    @{Sheet="WVR"; Cell="L132"; Value=17404.173},
    @{Sheet="WVR"; Cell="M132"; Value=-8676.049999999999},
    @{Sheet="WVR"; Cell="N132"; Value=-22464.173},
)

foreach ($u in $cellUpdates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

$cellDeletions = @(
    @{Sheet="LTW"; Cell="N59"},
)

foreach ($d in $cellDeletions) {
    $ws = $wb.Worksheets.Item($d.Sheet)
    $ws.Range($d.Cell).ClearContents()
}

Write-Host "Applied $(256) cell updates and $(1) cell deletions."
